# Insert a new "segments" column at B, pushing the existing metric columns
# (PercActivations .. totalStd) from B..K to C..L, then move the segment
# name strings that used to live in column A into the new column B and
# replace column A with a 0-based numeric segment index.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column before the old column B. Excel shifts the old
#    B:K header/data block to C:L for us (dimension becomes A1:L20).
$ws.Columns.Item(2).Insert()

# 2) The new B1 header cell needs the same (bold/border/centered) style
#    as the rest of the header row. Copy formats only from the
#    neighbouring header cell rather than setting font/border properties
#    directly, so this reuses the existing style instead of minting a
#    new one.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1").Value = "segments"

# 3) Move the segment-name strings out of column A into the new column B,
#    and put a 0-based numeric index into column A instead.
$names = @("background","back_bumper","back_glass","back_left_door","back_left_light","back_right_door","back_right_light","front_bumper","front_glass","front_left_door","front_left_light","front_right_door","front_right_light","hood","left_mirror","right_mirror","tailgate","trunk","wheel")

for ($i = 0; $i -lt $names.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 2).Value = $names[$i]
    $ws.Cells.Item($r, 1).Value = $i
}

# 4) Column B (segment names) shouldn't carry the bordered/bold style that
#    Insert() copied over from column A - clear formatting there so only
#    column A (the index) keeps that style, matching the header styling
#    rules used elsewhere in the sheet.
$ws.Range("B2:B20").ClearFormats()
